$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily")
$ws.Range("A19").Value = 44544
